$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table (rows 2-51) with freshly scraped values.
# Cells are written with a leading apostrophe to force text storage (so values like
# "30.734.49", "95.00" or "1.000" aren't auto-coerced into numbers and lose their
# exact textual formatting), then the style is reset to "Normal" so no stray
# quote-prefix formatting is left attached to the cell.
$ws.Range('D2').Value = "'30.734.49"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.91%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.876.76"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.14%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.00%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'237.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.98%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.04%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4785"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +2.35%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2834"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +4.89%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06504"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +3.89%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'18.79"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +17.26%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'1.889.74"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +2.83%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.36%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'95.00"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +13.64%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.102"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +3.74%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.6515"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +5.17%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'293.76"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +29.63%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'30.709.79"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.15%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = "'Dai"
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = "'1.001"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.09%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = "'Avalanche"
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = "'13.07"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +6.14%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.000007530"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +3.51%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'2.112.69"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.83%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.9991"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'5.148"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +5.91%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'6.094"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +4.48%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'169.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +3.14%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'9.230"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.43%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'19.53"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +9.94%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'1.969"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +4.71%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'0.1053"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +1.75%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'1.355"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -1.10%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'4.130"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.59%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'3.941"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +4.12%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.04971"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +3.63%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.176"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +3.63%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.7222"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.18%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.704"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.43%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.01933"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.98%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'2.715"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +2.66%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'2.058"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +7.57%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.8900"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.02%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'107.64"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +3.41%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'1.000"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.09%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.4203"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +5.27%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'5.576"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.09%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'7.335"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +4.93%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'65.41"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +9.92%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.1226"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +3.05%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'34.60"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +5.95%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'8.828"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +3.36%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'1.388"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +2.57%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.05559"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.98%  "
$ws.Range('E51').Style = 'Normal'
